# Generate Report for Handoff
# Updates the status of the "23a98bbb-da5a-43fe-856b-4adafeabc453" file
# from "Handed back: in sync with en-US" to "Ready for handoff", and
# records new handoff datetimes for zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the 23a98bbb file, columns B (zh-cn) and C (de-de)
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the 23a98bbb file
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-23 08:55:30"

# de-de sheet: row 3 is the 23a98bbb file
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-23 08:55:42"
